# Add a new item row ("MUCOPHYLLINE SYRUP 125 ML") to the pharmacy report,
# inserted in its alphabetically-correct position between "MILGA ADVANCE 30
# F.C. TABS" (row 13) and "PERLOC 40MG 14 F.C.TAB." (old row 14).
#
# This pushes the old rows 14-23 (items 11-18, the totals row and the
# footer row) down by one row, and the new row 14 carries the new item's
# data: balance "10:0", price 50, turnover count 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift rows 14 downward (this carries values, styles and merged cells
#    of every row from 14 to 23 down to 15-24 automatically).
$ws.Rows("14").Insert()

# 2) The newly-inserted row 14 is blank and does not inherit the table's
#    cell styles, so copy them over from row 15 (which now holds what used
#    to be row 14 - still formatted with the standard item-row style).
$ws.Range("A15:N15").Copy()
$ws.Range("A14:N14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Re-create the merged cells for the new row (Insert() does not clone
#    merge areas into the freshly-inserted row).
$ws.Range("B14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()

# 4) Match the row height used by the other item rows.
$ws.Rows("14").RowHeight = 25.5

# 5) Fill in the new item's data.
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "MUCOPHYLLINE SYRUP 125 ML"
$ws.Range("H14").Value = "10:0"
$ws.Range("L14").Value = 50
$ws.Range("N14").Value = 1

# 6) Update the running total (K column on the totals row, now row 23 after
#    the insert) to include the new item's price.
$ws.Range("K23").Value = $ws.Range("K23").Value + 50

# 7) Restore the row heights for the totals/footer rows that shifted down.
$ws.Rows("23").RowHeight = 25.5
$ws.Rows("24").RowHeight = 17.25
